# "adding term 2.0 now utf-8"
# Bumps the value set from 1.1.0 to 2.0.0: updates the metadata block on the
# "Metadata" sheet, and adds a new Concept row (with its generated concept id)
# on the "Include from FSIII" sheet, pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: bump version/date, fix the contact + description text
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B3").Value  = "2.0.0"
$wsMeta.Range("B8").Value  = "2024-06-03T10:45:43+02:00"
$wsMeta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"
$wsMeta.Range("B11").Value = "Matter of interest values to support when no observations have been made"

# ---------------------------------------------------------------------------
# "Include from FSIII" sheet: insert a new Concept ("d7ff926a-...") above the
# existing "B6" concept row, pushing the other rows (including "System URI")
# down by one.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Include from FSIII")

# Shift rows 2 downward, opening up a blank row 2 for the new concept (this
# keeps every other row - "B6", the blank description row, "System URI" -
# intact, just moved down one row).
$ws2.Range("A2:B2").Insert()

# Excel copies row 1's (header) formatting into the freshly inserted row;
# re-apply the plain body styling used by the rest of the table by pulling
# it back from the row directly below (now the "B6" row).
$ws2.Range("A3:B3").Copy($ws2.Range("A2:B2"))

# Row 2 now becomes the new concept entry; its Description cell stays blank.
$ws2.Range("A2").Value = "d7ff926a-4955-478f-b300-0b0ec0785013"
$ws2.Range("B2").ClearContents()

Write-Host "done"
